# Apply the commit's changes:
# - Typography sheet: fill Fallback Character / Wildcard Characters /
#   Wildcard Ranges example columns (F,G,H) for rows 4-6.
# - Translation sheet: add a queue of new text ids (rows 6-12) used for
#   the main screen settings (hour/minute, am/pm, day of week, etc.)

$wb = $excel.ActiveWorkbook

$wsTypo = $wb.Worksheets.Item("Typography")

$wsTypo.Range("F4").Value = "?"
$wsTypo.Range("G4").Value = "abcdefghijklmnopqrstuvwxyzABCDEFGHIJKLMNOPQRSTUVWXYZ0123456789"
$wsTypo.Range("H4").Value = "0-9,A-Z,a-z"
$wsTypo.Range("F4:H4").Style = "Normal"

$wsTypo.Range("F5").Value = "?"
$wsTypo.Range("G5").Value = "abcdefghijklmnopqrstuvwxyzABCDEFGHIJKLMNOPQRSTUVWXYZ0123456789"
$wsTypo.Range("H5").Value = "0-9,A-Z,a-z"
$wsTypo.Range("F5:H5").Style = "Normal"

$wsTypo.Range("F6").Value = "?"
$wsTypo.Range("G6").Value = "abcdefghijklmnopqrstuvwxyzABCDEFGHIJKLMNOPQRSTUVWXYZ0123456789"
$wsTypo.Range("H6").Value = "0-9,A-Z,a-z"
$wsTypo.Range("F6:H6").Style = "Normal"

$wsTrans = $wb.Worksheets.Item("Translation")

$wsTrans.Range("B6").Value = "SingleUseId3"
$wsTrans.Range("C6").Value = "Large"
$wsTrans.Range("D6").Value = "Right"
$wsTrans.Range("E6").Value = "<hour>:<minute>"
$wsTrans.Range("F6").Value = "LTR"
$wsTrans.Range("B6:F6").Style = "Normal"

$wsTrans.Range("B7").Value = "SingleUseId4"
$wsTrans.Range("C7").Value = "Large"
$wsTrans.Range("D7").Value = "Left"
$wsTrans.Range("E7").NumberFormat = "@"
$wsTrans.Range("E7").Value = "0"
$wsTrans.Range("F7").Value = "LTR"
$wsTrans.Range("B7:F7").Style = "Normal"

$wsTrans.Range("B8").Value = "SingleUseId5"
$wsTrans.Range("C8").Value = "Large"
$wsTrans.Range("D8").Value = "Left"
$wsTrans.Range("E8").NumberFormat = "@"
$wsTrans.Range("E8").Value = "00"
$wsTrans.Range("F8").Value = "LTR"
$wsTrans.Range("B8:F8").Style = "Normal"

$wsTrans.Range("B9").Value = "SingleUseId6"
$wsTrans.Range("C9").Value = "Large"
$wsTrans.Range("D9").Value = "Center"
$wsTrans.Range("E9").Value = "<hF>"
$wsTrans.Range("F9").Value = "LTR"
$wsTrans.Range("B9:F9").Style = "Normal"

$wsTrans.Range("B10").Value = "SingleUseId7"
$wsTrans.Range("C10").Value = "Large"
$wsTrans.Range("D10").Value = "Left"
$wsTrans.Range("E10").Value = "pm"
$wsTrans.Range("F10").Value = "LTR"
$wsTrans.Range("B10:F10").Style = "Normal"

$wsTrans.Range("B11").Value = "SingleUseId8"
$wsTrans.Range("C11").Value = "Default"
$wsTrans.Range("D11").Value = "Center"
$wsTrans.Range("E11").Value = "<dow>"
$wsTrans.Range("F11").Value = "LTR"
$wsTrans.Range("B11:F11").Style = "Normal"

$wsTrans.Range("B12").Value = "SingleUseId9"
$wsTrans.Range("C12").Value = "Default"
$wsTrans.Range("D12").Value = "Left"
$wsTrans.Range("E12").Value = "Monday"
$wsTrans.Range("F12").Value = "LTR"
$wsTrans.Range("B12:F12").Style = "Normal"
